$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "falls" / 폭포 quiz entry entirely (was row 4) — shifts all
# subsequent rows up by one.
$ws.Rows(4).Delete()

# The former "mushroom" / 버섯 entry (now row 9 after the shift above)
# becomes "bush" / 덤불.
$ws.Range("B9").Value = "bush"
$ws.Range("D9").Value = "bush"
$ws.Range("F9").Value = "덤불"

# Match the author's final selection.
[void]$ws.Range("F10").Select()
